# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while guaranteeing it stays plain text (no
# auto-coercion of numeric-looking strings like "586.80" -> 586.8),
# then drop the temporary Text number-format so the cells style
# index is unchanged from the original (General/default).
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '61.262.92'
$ws.Range('E2').Value = '  -2.58%  '
Set-TextValue 'D3' '2.961.70'
$ws.Range('E3').Value = '  -2.50%  '
$ws.Range('E4').Value = '  +0.10%  '
Set-TextValue 'D5' '586.80'
$ws.Range('E5').Value = '  +1.21%  '
Set-TextValue 'D6' '140.88'
$ws.Range('E6').Value = '  -6.28%  '
$ws.Range('E7').Value = '  +0.07%  '
Set-TextValue 'D8' '0.517'
$ws.Range('E8').Value = '  -2.10%  '
Set-TextValue 'D9' '2.957.99'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('E10').Value = '  -5.97%  '
Set-TextValue 'D11' '5.73'
$ws.Range('E11').Value = '  -0.47%  '
Set-TextValue 'D12' '0.454'
$ws.Range('E12').Value = '  +2.15%  '
Set-TextValue 'D13' '0.0000224'
$ws.Range('E13').Value = '  -3.28%  '
Set-TextValue 'D14' '33.75'
$ws.Range('E14').Value = '  -5.40%  '
$ws.Range('E15').Value = '  +1.33%  '
Set-TextValue 'D16' '3.461.51'
$ws.Range('E16').Value = '  -2.35%  '
Set-TextValue 'D17' '6.96'
$ws.Range('E17').Value = '  -1.59%  '
Set-TextValue 'D18' '61.394.36'
$ws.Range('E18').Value = '  -2.40%  '
Set-TextValue 'D19' '2.969.25'
$ws.Range('E19').Value = '  -2.51%  '
Set-TextValue 'D20' '446.98'
$ws.Range('E20').Value = '  -5.93%  '
Set-TextValue 'D21' '13.78'
$ws.Range('E21').Value = '  -2.71%  '
Set-TextValue 'D22' '0.679'
$ws.Range('E22').Value = '  -2.83%  '
Set-TextValue 'D23' '7.29'
$ws.Range('E23').Value = '  -2.08%  '
Set-TextValue 'D24' '80.97'
$ws.Range('E24').Value = '  -0.25%  '
Set-TextValue 'D25' '12.05'
$ws.Range('E25').Value = '  -3.72%  '
$ws.Range('E26').Value = '  -9.23%  '
Set-TextValue 'D27' '0.999'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  -6.71%  '
$ws.Range('E29').Value = '  -0.06%  '
Set-TextValue 'D30' '2.62'
$ws.Range('E30').Value = '  -0.22%  '
Set-TextValue 'D31' '6.81'
$ws.Range('E31').Value = '  -6.22%  '
$ws.Range('E32').Value = '  -6.09%  '
Set-TextValue 'D33' '26.82'
$ws.Range('E33').Value = '  -2.33%  '
$ws.Range('E34').Value = '  -3.01%  '
$ws.Range('E35').Value = '  -3.77%  '
Set-TextValue 'D36' '0.0₃0766'
$ws.Range('E36').Value = '  -4.60%  '
Set-TextValue 'D37' '5.68'
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D38' '50.00'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '2.05'
$ws.Range('E39').Value = '  -5.15%  '
Set-TextValue 'D40' '9.07'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('E41').Value = '  +4.18%  '
Set-TextValue 'D42' '2.75'
$ws.Range('E42').Value = '  -10.71%  '
Set-TextValue 'D43' '385.69'
$ws.Range('E43').Value = '  -8.71%  '
$ws.Range('E44').Value = '  -1.88%  '
Set-TextValue 'D45' '2.706.09'
$ws.Range('E45').Value = '  -4.18%  '
Set-TextValue 'D46' '0.261'
$ws.Range('E46').Value = '  -7.62%  '
Set-TextValue 'D47' '36.65'
$ws.Range('E47').Value = '  -2.93%  '
Set-TextValue 'D48' '129.88'
$ws.Range('E48').Value = '  +2.55%  '
$ws.Range('E50').Value = '  -1.44%  '
Set-TextValue 'D51' '2.14'
$ws.Range('E51').Value = '  -1.18%  '
